# Remove the stray duplicate "Rectangle 73" placeholder shapes (the ones with
# the hard-coded / stale "xx" / "XXX" / "202"+"6" year text and the
# x=1218,y=-5264 offset) from each slide that has one. These are leftover
# duplicates of the real "top-left year" shape; removing them lets the
# remaining shape update the year automatically instead of having two shapes
# fight over the same spot.
#
# Target shape ids (unique within their slide, per PowerPoint's Shape.Id):
#   Slide 2 -> id 119
#   Slide 4 -> id 19
#   Slide 5 -> id 61
#   Slide 6 -> id 107
#   Slide 7 -> id 110
# Note: Slide 2 also has another shape named "Rectangle 73" (id 95) that is a
# different shape (different position/content) and must be left untouched.

$p = $ppt.ActivePresentation

$targets = @{
    2 = 119
    4 = 19
    5 = 61
    6 = 107
    7 = 110
}

foreach ($slideIndex in $targets.Keys) {
    $s = $p.Slides.Item($slideIndex)
    $targetId = $targets[$slideIndex]

    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($sh.Name -eq "Rectangle 73" -and $sh.Id -eq $targetId) {
            $sh.Delete()
            break
        }
    }
}
